$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 3134
$ws1.Range("F8").Value = 2724
$ws1.Range("F13").Value = 279
$ws1.Range("F15").Value = 5594
$ws1.Range("F19").Value = 78
$ws1.Range("F24").Value = 325

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 997

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F9").Value = 1424

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1424
$ws4.Range("F15").Value = 3134
$ws4.Range("F16").Value = 2724
$ws4.Range("F22").Value = 279
$ws4.Range("F24").Value = 5594
$ws4.Range("F30").Value = 78
$ws4.Range("F48").Value = 325
